$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Date Deployed" cell (A1) with the new date
$ws.Range("A1").Value = "Date Deployed: 31/3/2019"

# Reset the view: scroll back to top-left and select A2
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
